$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("BCS")
$ws3.Name = "BCS-BCS"
$newSheet = $wb.Worksheets.Add($null, $ws3)
$newSheet.Name = "BCS-DoSfCS"

$newSheet.Range("A1").Value = "Years"
$years = 2021..2050
for ($i = 0; $i -lt $years.Count; $i++) {
    $col = $i + 2
    $newSheet.Cells.Item(1, $col).Value = $years[$i]
}

$newSheet.Range("A2").Value = "Duration"
$newSheet.Range("B2").Value = 12
$newSheet.Range("C2").Formula = "=`$B`$2"
$newSheet.Range("D2").Formula = "=`$B`$2"
$newSheet.Range("E2:AE2").Formula = "=`$B`$2"
